$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "SCORE" header column (H2), matching the formatting of the
# existing header cells by copying G2's format, then replacing its value.
$ws.Range("G2").Copy($ws.Range("H2"))
$ws.Range("H2").Value = "SCORE"

# Update the selection to reflect the new active cell after the edit
$ws.Range("I5").Select()
